# Automatische test-sync: 2025-08-28 21:01:50
# Appends a new log entry (row 18) to the "Logs" sheet, extends the
# conditional-formatting ranges that covered rows 2-17 so they also cover
# row 18, and bumps the "Retour / Terugbetaling" tally on the Dashboard
# sheet from 15 to 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 18

$ws.Cells.Item($newRow, 1).Value  = "Retour status"
$ws.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 4).Value  = "Retour / Terugbetaling"
$ws.Cells.Item($newRow, 6).Value  = "2025-08-28 21:01:03"
$ws.Cells.Item($newRow, 7).Value  = "Nee"
$ws.Cells.Item($newRow, 8).Value  = "Ja"
$ws.Cells.Item($newRow, 9).Value  = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional formatting that applied to rows 2:17 on columns
# D, G, H, I and J so it now also covers the freshly added row 18.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "17")
    $newRange = $ws.Range($col + "2:" + $col + "18")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Keep the Dashboard summary count in sync with the new log entry.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 16
